# Commit: "added folders for the processing"
# Breda.docx previously listed every extracted measure-category bullet point
# one paragraph at a time; this collapses each of the four answer blocks into
# a single narrative paragraph explaining why "Geen" (none) applies, and
# removes the now-redundant bullet list in between.

$d = $word.ActiveDocument

# Locate the target paragraphs by their current text (paragraph mark trimmed)
# rather than hard-coded indices, so the script is resilient to how the
# runtime numbers things.
function Get-ParagraphIndexByText($doc, $targetText, $occurrence) {
    $seen = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            $seen++
            if ($seen -eq $occurrence) {
                return $i
            }
        }
    }
    return -1
}

# Resolve every target index up front (before any mutation), so earlier
# edits never shift the meaning of a later lookup.
$idx1 = Get-ParagraphIndexByText $d "Geen." 1
$idx2 = Get-ParagraphIndexByText $d "Geen." 2
$idx3 = Get-ParagraphIndexByText $d "- Mobiliteit (verkeer)" 1
$lastIdx = $d.Paragraphs.Count

# 1) First "Geen." paragraph -> long explanation of why no air-quality measures were found
$d.Paragraphs.Item($idx1).Range.Text = 'Er zijn geen maatregelen gerelateerd aan het verbeteren van de luchtkwaliteit in deze tekst. De categorie voor alle maatregelen is daarom "Geen". Hoewel er wel wordt gesproken over klimaatadaptatie en het belang van groen en natuur voor de gezondheid, worden er geen concrete maatregelen genoemd om de luchtkwaliteit te verbeteren.'

# 2) Second "Geen." paragraph -> its own explanation
$d.Paragraphs.Item($idx2).Range.Text = 'Er zijn geen maatregelen genoemd die specifiek gericht zijn op het verbeteren van de luchtkwaliteit.'

# 3) First "- Mobiliteit (verkeer)" bullet (start of the long bullet list) -> explanation
$d.Paragraphs.Item($idx3).Range.Text = 'Er worden geen specifieke maatregelen genoemd die gericht zijn op het verbeteren van de luchtkwaliteit in deze tekst. De tekst richt zich voornamelijk op stedelijke ontwikkelingsprojecten en -plannen, mobiliteit, burger- en bedrijfsparticipatie, monitoring, locaties met hoge blootstelling en kwetsbare groepen, en internationaal luchtbeleid. Er worden wel maatregelen genoemd om de mobiliteit te verbeteren, zoals het elektrificeren van vervoer en het beperken van effecten op de omgeving. Daarnaast worden er maatregelen genoemd om geurhinder te verminderen en woningbouw en bedrijfsactiviteiten op voldoende afstand van elkaar te plannen om geurhinder te voorkomen in de industrie en het uitsluiten van risicobedrijven op nieuwe bedrijfsterreinen en het clusteren van grote industriële risicobedrijven op industrieterrein Moerdijk.'

# 4) Delete the rest of the bullet list: everything between the paragraph we
#    just rewrote and the final summary paragraph (which is rewritten next).
$startPar = $d.Paragraphs.Item($idx3 + 1)
$endPar = $d.Paragraphs.Item($lastIdx - 1)
$rangeToDelete = $d.Range($startPar.Range.Start, $endPar.Range.End)
$rangeToDelete.Delete()

# 5) Final summary paragraph (now the last paragraph in the document) -> explanation
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = 'Er worden geen maatregelen genoemd die gerelateerd zijn aan het verbeteren van de luchtkwaliteit in de tekst. Daarom is het antwoord "Geen". Mobiliteit (verkeer) omvat onder andere het verbeteren van de toegang tot het platteland, het creëren van groene routes, mobiliteitshubs en het stimuleren van initiatieven van burgers. Voor mobiele machines, industrie, houtverbranding in particuliere huishoudens, binnenvaart en havens, landbouw, monitoring en kwetsbare groepen worden geen specifieke maatregelen genoemd. Burger- en bedrijfsparticipatie omvat onder andere het verbeteren van de dienstverlening van de gemeente en het betrekken van burgers bij het beheer van de openbare ruimte. Het internationaal luchtkwaliteitsbeleid wordt niet besproken in de tekst.'

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
